$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Update cached date field text ("10.12.2013" -> "11.05.2014")
#    on the slide master and all 11 slide layouts (footer "dt"
#    placeholder, PlaceholderFormat.Type = 16 = ppPlaceholderDate).
# ---------------------------------------------------------------
$master = $p.SlideMaster

$phCount = $master.Shapes.Placeholders.Count
for ($i = 1; $i -le $phCount; $i++) {
    $ph = $master.Shapes.Placeholders.Item($i)
    if ($ph.PlaceholderFormat.Type -eq 16) {
        $ph.TextFrame.TextRange.Text = "11.05.2014"
    }
}

$layoutCount = $master.CustomLayouts.Count
for ($li = 1; $li -le $layoutCount; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $phCount2 = $layout.Shapes.Placeholders.Count
    for ($i = 1; $i -le $phCount2; $i++) {
        $ph = $layout.Shapes.Placeholders.Item($i)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $ph.TextFrame.TextRange.Text = "11.05.2014"
        }
    }
}

Write-Host "Dates updated"

# ---------------------------------------------------------------
# 2. Slide 1: subtitle placeholder - reposition/resize, enable
#    shrink-text-on-overflow autofit, and rewrite the author list
#    as four right-aligned paragraphs.
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)

$subtitle.Left = 108.0
$subtitle.Top = 287.00976377952753
$subtitle.Width = 504.0
$subtitle.Height = 185.12693033385827
$subtitle.TextFrame.AutoSize = 2
$subtitle.TextFrame.TextRange.Text = "Демьяненко Илья`rПопов Кирилл`rСоболев Артем`r444 группа"

Write-Host "Slide 1 subtitle updated"

# ---------------------------------------------------------------
# 3. Slide 2: "Vkontakte" -> "ВКонтакте" inside the API paragraph.
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$content2 = $s2.Shapes.Item(2)
$tr2 = $content2.TextFrame.TextRange
$target2 = $tr2.Characters(108, 9)
$target2.Text = "ВКонтакте"

Write-Host "Slide 2 updated"

# ---------------------------------------------------------------
# 4. Slide 6: title "Бизнес функции" -> "Бизнес-функции".
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$title6 = $s6.Shapes.Item(1)
$title6.TextFrame.TextRange.Text = "Бизнес-функции"

Write-Host "Slide 6 updated"
